$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.308.39"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  -3.03%  "

$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.005.64"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -3.45%  "

$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  +0.02%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.36"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -1.88%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.88"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -6.58%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -3.41%  "

$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.003.40"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -3.54%  "

$ws.Range("E10").Value = "  -6.25%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.65"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -4.47%  "

$ws.Range("E12").Value = "  -2.58%  "

$ws.Range("E13").Value = "  -5.16%  "

$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.60"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -7.16%  "

$ws.Range("E15").Value = "  +1.66%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.501.89"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -3.31%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.10"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -1.98%  "

$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.329.86"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -2.72%  "

$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.005.43"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -3.23%  "

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "454.06"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -5.05%  "

$ws.Range("E21").Value = "  -4.32%  "

$ws.Range("E22").Value = "  -4.47%  "

$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.32"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -3.76%  "

$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.29"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -6.96%  "

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.12"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -1.30%  "

$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.31"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -5.05%  "

$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -3.39%  "

$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -0.17%  "

$ws.Range("E29").Value = "  +0.10%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.17"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -4.47%  "

$ws.Range("E31").Value = "  -2.76%  "

$ws.Range("E32").Value = "  -4.74%  "

$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.90"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -1.95%  "

$ws.Range("E34").Value = "  -5.19%  "

$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.03"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -3.13%  "

$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0791"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -6.41%  "

$ws.Range("E37").Value = "  -4.93%  "

$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.12"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -5.27%  "

$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.13"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -1.85%  "

$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.06"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -1.91%  "

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.93"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -11.30%  "

$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "409.73"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -6.89%  "

$ws.Range("E43").Value = "  -5.55%  "

$ws.Range("E44").Value = "  -1.25%  "

$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.771.09"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -2.09%  "

$ws.Range("E46").Value = "  -3.68%  "

$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.85"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  -6.31%  "

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.11"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -1.52%  "

$ws.Range("E50").Value = "  -2.07%  "

$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.82"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -6.55%  "
